$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on price cells whose new values would otherwise be
# auto-converted to numbers by Excel (losing formatting like trailing zeros).
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Update Price (D) and Volume(1h) (E) columns per latest crypto snapshot.
$ws.Range("D2").Value = "22.376.37"
$ws.Range("E2").Value = "  -4.19%  "
$ws.Range("D3").Value = "1.569.44"
$ws.Range("E3").Value = "  -3.57%  "
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("E5").Value = "  -0.05%  "
$ws.Range("D6").Value = "289.00"
$ws.Range("E6").Value = "  -2.98%  "
$ws.Range("D7").Value = "0.3671"
$ws.Range("E7").Value = "  -2.10%  "
$ws.Range("D8").Value = "49.26"
$ws.Range("E8").Value = "  -2.44%  "
$ws.Range("D9").Value = "0.3385"
$ws.Range("E9").Value = "  -2.44%  "
$ws.Range("D10").Value = "1.168"
$ws.Range("E10").Value = "  -2.36%  "
$ws.Range("D11").Value = "0.07622"
$ws.Range("E11").Value = "  -4.90%  "
$ws.Range("E12").Value = "  -0.04%  "
$ws.Range("D13").Value = "21.22"
$ws.Range("E13").Value = "  -2.67%  "
$ws.Range("D14").Value = "6.057"
$ws.Range("E14").Value = "  -3.57%  "
$ws.Range("D15").Value = "6.908"
$ws.Range("E15").Value = "  -3.99%  "
$ws.Range("D16").Value = "1.576.54"
$ws.Range("E16").Value = "  -3.22%  "
$ws.Range("D17").Value = "0.00001131"
$ws.Range("E17").Value = "  -4.31%  "
$ws.Range("D18").Value = "89.58"
$ws.Range("E18").Value = "  -5.07%  "
$ws.Range("D19").Value = "0.06751"
$ws.Range("E19").Value = "  -2.63%  "
$ws.Range("E20").Value = "  -0.08%  "
$ws.Range("D21").Value = "6.219"
$ws.Range("E21").Value = "  -5.59%  "
$ws.Range("D22").Value = "0.5321"
$ws.Range("E22").Value = "  -5.96%  "
$ws.Range("D23").Value = "16.51"
$ws.Range("E23").Value = "  -4.11%  "
$ws.Range("D24").Value = "11.97"
$ws.Range("E24").Value = "  -2.74%  "
$ws.Range("D25").Value = "22.375.79"
$ws.Range("E25").Value = "  -4.25%  "
$ws.Range("D26").Value = "2.368"
$ws.Range("E26").Value = "  -2.76%  "
$ws.Range("D27").Value = "2.894"
$ws.Range("E27").Value = "  -4.16%  "
$ws.Range("D28").Value = "20.03"
$ws.Range("E28").Value = "  -3.33%  "
$ws.Range("D29").Value = "145.16"
$ws.Range("E29").Value = "  -3.93%  "
$ws.Range("D30").Value = "4.975"
$ws.Range("E30").Value = "  -3.59%  "
$ws.Range("D31").Value = "125.39"
$ws.Range("E31").Value = "  -4.40%  "
$ws.Range("D32").Value = "1.745.08"
$ws.Range("E32").Value = "  -3.75%  "
$ws.Range("D33").Value = "1.038"
$ws.Range("E33").Value = "  +6.76%  "
$ws.Range("D34").Value = "6.237"
$ws.Range("E34").Value = "  -6.50%  "
$ws.Range("D35").Value = "2.017"
$ws.Range("E35").Value = "  -5.28%  "
$ws.Range("D36").Value = "10.19"
$ws.Range("E36").Value = "  -9.74%  "
$ws.Range("D37").Value = "0.08455"
$ws.Range("E37").Value = "  -2.92%  "
$ws.Range("D38").Value = "0.02534"
$ws.Range("E38").Value = "  -4.46%  "
$ws.Range("E39").Value = "  -3.85%  "
$ws.Range("D40").Value = "5.524"
$ws.Range("E40").Value = "  -4.88%  "
$ws.Range("D41").Value = "0.06475"
$ws.Range("E41").Value = "  -2.58%  "
$ws.Range("D42").Value = "1.314"
$ws.Range("E42").Value = "  +2.35%  "
$ws.Range("D43").Value = "11.70"
$ws.Range("E43").Value = "  -7.19%  "
$ws.Range("E44").Value = "  -6.22%  "
$ws.Range("D45").Value = "14.14"
$ws.Range("E45").Value = "  -7.75%  "
$ws.Range("E46").Value = "  +0.00%  "
$ws.Range("D47").Value = "0.5972"
$ws.Range("E47").Value = "  -4.91%  "
$ws.Range("D48").Value = "3.752"
$ws.Range("E48").Value = "  -3.50%  "
$ws.Range("D49").Value = "2.097"
$ws.Range("E49").Value = "  -5.63%  "
$ws.Range("E50").Value = "  +4.79%  "
$ws.Range("D51").Value = "124.68"
$ws.Range("E51").Value = "  -1.35%  "
